$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column for all existing data
#    rows (2-34): 45629 -> 45630 (one day later).
$ws.Range("C2:C34").Value2 = 45630

# 2. Row 34 picks up an explicit row height in the new file (it previously
#    relied on the sheet default). Setting it matches the diff exactly.
$ws.Rows.Item(34).RowHeight = 15

# 3. Append the new complaint row (row 35).
$ws.Range("A35").Value = "A 57397-2024"

$ws.Range("B35").Value2 = 45629
$ws.Range("B35").NumberFormat = "YYYY-MM-DD"

$ws.Range("C35").Value2 = 45630
$ws.Range("C35").NumberFormat = "YYYY-MM-DD"

$ws.Range("D35").Value = "OKÄNT"
$ws.Range("E35").Value = "OKÄNT"

$ws.Range("G35").Value2 = 3.8
$ws.Range("H35:Q35").Value2 = 0

# R35 stays an empty (wrap-text styled) cell, matching the style used for
# the same column in every other data row.
$ws.Range("R35").WrapText = $true
